$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 239.0839323333333
$ws.Range("H2").Value = 717.251797
$ws.Range("I2").Value = 0.4086975387666237
$ws.Range("J2").Value = 0.4086975387666237
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.490547666666667
$ws.Range("N2").Value = 19.471643
$ws.Range("O2").Value = 0.8021666724616637
$ws.Range("P2").Value = 0.8021666724616636
$ws.Range("Q2").Value = 1551.785659143608
$ws.Range("R2").Value = 13966.07093229247
$ws.Range("S2").Value = 0.3278435447156943
$ws.Range("T2").Value = 0.3278435447156943

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 239.0839323333333
$ws.Range("H3").Value = 717.251797
$ws.Range("I3").Value = 0.4086975387666237
$ws.Range("J3").Value = 0.4086975387666237
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.600723
$ws.Range("N3").Value = 4.802169
$ws.Range("O3").Value = 0.1978333275383364
$ws.Range("P3").Value = 0.1978333275383364
$ws.Range("Q3").Value = 382.7071494164103
$ws.Range("R3").Value = 3444.364344747693
$ws.Range("S3").Value = 0.0808539940509294
$ws.Range("T3").Value = 0.0808539940509294

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 117.0512696666667
$ws.Range("H4").Value = 351.153809
$ws.Range("I4").Value = 0.2000910950200451
$ws.Range("J4").Value = 0.2000910950200451
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.490547666666667
$ws.Range("N4").Value = 19.471643
$ws.Range("O4").Value = 0.8021666724616637
$ws.Range("P4").Value = 0.8021666724616636
$ws.Range("Q4").Value = 759.7268452153542
$ws.Range("R4").Value = 6837.541606938187
$ws.Range("S4").Value = 0.1605064078814401
$ws.Range("T4").Value = 0.1605064078814401

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 117.0512696666667
$ws.Range("H5").Value = 351.153809
$ws.Range("I5").Value = 0.2000910950200451
$ws.Range("J5").Value = 0.2000910950200451
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.600723
$ws.Range("N5").Value = 4.802169
$ws.Range("O5").Value = 0.1978333275383364
$ws.Range("P5").Value = 0.1978333275383364
$ws.Range("Q5").Value = 187.3666595346357
$ws.Range("R5").Value = 1686.299935811721
$ws.Range("S5").Value = 0.03958468713860497
$ws.Range("T5").Value = 0.03958468713860496

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 171.15883
$ws.Range("H6").Value = 513.47649
$ws.Range("I6").Value = 0.2925842480357353
$ws.Range("J6").Value = 0.2925842480357353
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 6.490547666666667
$ws.Range("N6").Value = 19.471643
$ws.Range("O6").Value = 0.8021666724616637
$ws.Range("P6").Value = 0.8021666724616636
$ws.Range("Q6").Value = 1110.914544685897
$ws.Range("R6").Value = 9998.23090217307
$ws.Range("S6").Value = 0.2347013326615238
$ws.Range("T6").Value = 0.2347013326615238

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 171.15883
$ws.Range("H7").Value = 513.47649
$ws.Range("I7").Value = 0.2925842480357353
$ws.Range("J7").Value = 0.2925842480357353
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.600723
$ws.Range("N7").Value = 4.802169
$ws.Range("O7").Value = 0.1978333275383364
$ws.Range("P7").Value = 0.1978333275383364
$ws.Range("Q7").Value = 273.97787583409
$ws.Range("R7").Value = 2465.80088250681
$ws.Range("S7").Value = 0.05788291537421148
$ws.Range("T7").Value = 0.05788291537421147

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 57.695868
$ws.Range("H8").Value = 173.087604
$ws.Range("I8").Value = 0.09862711817759588
$ws.Range("J8").Value = 0.09862711817759588
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.490547666666667
$ws.Range("N8").Value = 19.471643
$ws.Range("O8").Value = 0.8021666724616637
$ws.Range("P8").Value = 0.8021666724616636
$ws.Range("Q8").Value = 374.477781423708
$ws.Range("R8").Value = 3370.300032813372
$ws.Range("S8").Value = 0.07911538720300536
$ws.Range("T8").Value = 0.07911538720300534

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 57.695868
$ws.Range("H9").Value = 173.087604
$ws.Range("I9").Value = 0.09862711817759588
$ws.Range("J9").Value = 0.09862711817759588
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.600723
$ws.Range("N9").Value = 4.802169
$ws.Range("O9").Value = 0.1978333275383364
$ws.Range("P9").Value = 0.1978333275383364
$ws.Range("Q9").Value = 92.355102912564
$ws.Range("R9").Value = 831.1959262130761
$ws.Range("S9").Value = 0.01951173097459054
$ws.Range("T9").Value = 0.01951173097459053
